$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 486
$wsExpo.Range("F3").Value = 5786
$wsExpo.Range("F4").Value = 390
$wsExpo.Range("F5").Value = 75
$wsExpo.Range("F6").Value = 99
$wsExpo.Range("F9").Value = 545

# Sheet "全部类型" (all types) - same events duplicated, same column F updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 486
$wsAll.Range("F3").Value = 5786
$wsAll.Range("F4").Value = 390
$wsAll.Range("F6").Value = 75
$wsAll.Range("F7").Value = 99
$wsAll.Range("F11").Value = 545
